# Append the next 3 days of GSC export data to the "Chart" sheet.
# Dates are stored as literal text (matching the existing A-column cells),
# so a leading apostrophe is used to stop Excel from auto-converting the
# text into a date serial number; ClearFormats() then drops the
# "quote prefix" formatting that the apostrophe entry would otherwise add,
# leaving the cell on the sheet's default (General) style like its
# neighbours.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$ws.Cells.Item(45, 1).Formula = "'2025-11-18"
$ws.Cells.Item(45, 1).ClearFormats()
$ws.Cells.Item(45, 2).Value = 0
$ws.Cells.Item(45, 3).Value = 26

$ws.Cells.Item(46, 1).Formula = "'2025-11-19"
$ws.Cells.Item(46, 1).ClearFormats()
$ws.Cells.Item(46, 2).Value = 0
$ws.Cells.Item(46, 3).Value = 26

$ws.Cells.Item(47, 1).Formula = "'2025-11-20"
$ws.Cells.Item(47, 1).ClearFormats()
$ws.Cells.Item(47, 2).Value = 0
$ws.Cells.Item(47, 3).Value = 25
